$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.951.82"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "2.257.87"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.11"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.61"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  -1.08%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.16"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0831"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "2.604.23"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.46"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "2.258.71"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "43.891.34"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("E19").Value = "  -7.07%  "
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.55"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.77"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.13"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.25"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.09"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.17"
$ws.Range("E29").Value = "  +3.89%  "
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.49"
$ws.Range("E31").Value = "  +3.90%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.18"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0852"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.71"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("E35").Value = "  +11.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.95"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.05"
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.29"
$ws.Range("E39").Value = "  +20.70%  "
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "1.814.14"
$ws.Range("E44").Value = "  +4.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "75.84"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.38"
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.00"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.44"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  +6.35%  "
